$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.862.08'
$ws.Range('E2').Value = '  -2.43%  '

$ws.Range('D3').Value = '1.808.05'
$ws.Range('E3').Value = '  -3.05%  '

$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('D5').Value = "'230.80"
$ws.Range('E5').Value = '  -0.29%  '

$ws.Range('D6').Value = "'0.605"
$ws.Range('E6').Value = '  -1.43%  '

$ws.Range('E7').Value = '  +0.21%  '

$ws.Range('D8').Value = "'39.08"
$ws.Range('E8').Value = '  -8.72%  '

$ws.Range('E9').Value = '  +2.39%  '

$ws.Range('E10').Value = '  -2.83%  '

$ws.Range('E11').Value = '  -2.01%  '

$ws.Range('D12').Value = '2.070.39'
$ws.Range('E12').Value = '  -3.10%  '

$ws.Range('D13').Value = '1.802.88'
$ws.Range('E13').Value = '  -3.68%  '

$ws.Range('E14').Value = '  -3.32%  '

$ws.Range('D15').Value = "'10.86"
$ws.Range('E15').Value = '  -7.39%  '

$ws.Range('E16').Value = '  -4.50%  '

$ws.Range('D17').Value = '34.840.17'
$ws.Range('E17').Value = '  -2.53%  '

$ws.Range('D18').Value = "'69.27"
$ws.Range('E18').Value = '  -1.97%  '

$ws.Range('D19').Value = '0.0₃0780'
$ws.Range('E19').Value = '  -3.40%  '

$ws.Range('D20').Value = "'239.17"
$ws.Range('E20').Value = '  -4.13%  '

$ws.Range('D21').Value = "'11.74"
$ws.Range('E21').Value = '  -4.90%  '

$ws.Range('D22').Value = "'4.64"
$ws.Range('E22').Value = '  -2.65%  '

$ws.Range('E23').Value = '  +0.23%  '

$ws.Range('E24').Value = '  -0.81%  '

$ws.Range('D25').Value = "'173.63"
$ws.Range('E25').Value = '  +1.57%  '

$ws.Range('D26').Value = "'7.76"
$ws.Range('E26').Value = '  -3.36%  '

$ws.Range('D27').Value = "'17.18"
$ws.Range('E27').Value = '  -4.15%  '

$ws.Range('E28').Value = '  -3.46%  '

$ws.Range('D29').Value = "'1.51"
$ws.Range('E29').Value = '  +4.75%  '

$ws.Range('E30').Value = '  +0.23%  '

$ws.Range('E31').Value = '  +0.09%  '

$ws.Range('D32').Value = "'0.0545"
$ws.Range('E32').Value = '  -0.72%  '

$ws.Range('E33').Value = '  -4.41%  '

$ws.Range('B34').Value = 'TrustWalletToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D34').Value = "'1.17"
$ws.Range('E34').Value = '  +7.11%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.75"
$ws.Range('E35').Value = '  -8.37%  '

$ws.Range('D36').Value = "'0.682"
$ws.Range('E36').Value = '  -1.52%  '

$ws.Range('D37').Value = "'90.64"
$ws.Range('E37').Value = '  -9.96%  '

$ws.Range('D38').Value = "'1.34"
$ws.Range('E38').Value = '  +6.32%  '

$ws.Range('D39').Value = '1.311.15'
$ws.Range('E39').Value = '  -4.09%  '

$ws.Range('E40').Value = '  -3.35%  '

$ws.Range('E41').Value = '  -0.99%  '

$ws.Range('D42').Value = "'0.955"
$ws.Range('E42').Value = '  -6.10%  '

$ws.Range('D43').Value = "'14.08"
$ws.Range('E43').Value = '  -6.29%  '

$ws.Range('D44').Value = "'2.18"
$ws.Range('E44').Value = '  -13.04%  '

$ws.Range('E45').Value = '  -4.87%  '

$ws.Range('E46').Value = '  -2.92%  '

$ws.Range('E47').Value = '  -1.91%  '

$ws.Range('D48').Value = '1.991.19'
$ws.Range('E48').Value = '  -2.21%  '

$ws.Range('B49').Value = 'PaxDollar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D49').Value = "'1.01"
$ws.Range('E49').Value = '  +0.19%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0671"
$ws.Range('E50').Value = '  +7.07%  '

$ws.Range('D51').Value = "'98.66"
$ws.Range('E51').Value = '  -6.23%  '

# Reset style (remove quotePrefix formatting) to match original (unstyled) cells
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
